# Update the "Förändrad" (Changed) date column (C) for rows 2-28 from
# 45550 (2024-09-15) to 45551 (2024-09-16), keeping existing cell styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45550) {
        $cell.Value2 = 45551
    }
}
